# complianceReport.xlsx - add table body builder
# -----------------------------------------------------------------------
# The existing report header used a blank spacer row (row 1) and pushed the
# "license / issue-date / reporting-period / cycle-total" block out to
# columns I:L. This edit:
#   1. removes the blank spacer row so the title block starts at row 1,
#   2. pulls the right-hand info block in from I:L to G:J,
#   3. re-labels "Evan Hiner" as "Hiner, Evan",
#   4. turns the "Label: value" cells into bold-label / normal-value text,
#   5. formats the title / name / cycle-total blocks with the right
#      fonts & alignment, and
#   6. builds the compliance table header row (DATE / TITLE / SPONSOR /
#      DELIVERY METHOD / GENERAL / ETHICS STATE).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. drop the blank spacer row, shifting everything up one row -------
$ws.Rows("1").Delete()

# --- 2. drop the two spacer columns (G:H) so I:L becomes G:J ------------
$ws.Columns("G:H").Delete()

# --- 3. rename the practitioner cell -------------------------------------
$ws.Range("A3").Value = "Hiner, Evan"

# --- 4. rebuild the "Label: value" cells as bold-label/plain-value text -
function Set-LabelValue {
    param($cell, [string]$label, [string]$value)
    $full = $label + $value
    $cell.Value = $full
    $cell.Characters(1, $label.Length).Font.Bold = $true
}

Set-LabelValue -cell $ws.Range("G3") -label "License #: " -value "123123123"
Set-LabelValue -cell $ws.Range("G4") -label "Issue Date: " -value "12/31/2017"
Set-LabelValue -cell $ws.Range("G5") -label "Reporting Period: " -value "12/30/2016 - 12/31/2017"
Set-LabelValue -cell $ws.Range("A6") -label "Cycle Total: " -value "12/30/2016 - 12/31/2017"

# --- 5. formatting for the four header blocks ----------------------------
# "New York CPA" title block
$title = $ws.Range("A1:D2")
$title.Font.Bold = $true
$title.Font.Size = 12
$title.HorizontalAlignment = -4131   # xlLeft
$title.VerticalAlignment = -4108     # xlCenter

# "Page 1" / license-info column
$infoCells = @($ws.Range("J1"), $ws.Range("G3:J3"), $ws.Range("G4:J4"), $ws.Range("G5:J5"))
foreach ($c in $infoCells) {
    $c.HorizontalAlignment = -4152   # xlRight
    $c.VerticalAlignment = -4108     # xlCenter
}

# practitioner name block
$name = $ws.Range("A3:F5")
$name.Font.Bold = $true
$name.Font.Size = 16
$name.HorizontalAlignment = -4131   # xlLeft
$name.VerticalAlignment = -4108     # xlCenter

# "Cycle Total" block
$cycle = $ws.Range("A6:D7")
$cycle.HorizontalAlignment = -4131  # xlLeft
$cycle.VerticalAlignment = -4108    # xlCenter

# --- 6. compliance table header row (row 8) ------------------------------
$ws.Range("A8").Value = "DATE"
$ws.Range("B8").Value = "TITLE"
$ws.Range("C8").Value = "SPONSOR"
$ws.Range("D8").Value = "DELIVERY METHOD"
$ws.Range("E8").Value = "GENERAL"
$ws.Range("F8").Value = "ETHICS STATE"

# --- 7. page setup: narrower L/R margins, keep "fit to one page" --------
$ws.PageSetup.LeftMargin = 0.25 * 72
$ws.PageSetup.RightMargin = 0.25 * 72
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

Write-Host "edit complete"
